$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Steps" value for the existing last row (row 69)
$ws.Range("E69").Value = "C++: vector, std::sort, string.fin(c) != std::string::npos"

# Add new row 70 (value for A70/B70 first)
$ws.Range("A70").Value = 1491
$ws.Range("B70").Value = "Average Salary Excluding the Min and Max Salaray"

# Add new row 71
$ws.Range("A71").Value = 1733
$ws.Range("B71").Value = "Minimum Number of People to Teach"

# C70 is filled in last so its shared-string index comes after B71's
$ws.Range("C70").Value = "3 Pointers/Math"

# Match final selection/active cell shown in the diff
$ws.Range("C70").Select()
